# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data updates to the Diabolos_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H96").Value = 916
$ws.Range("I96").Value = 599.2
$ws.Range("K96").Value = 1797.6
$ws.Range("M96").Value = -424.6000000000001
$ws.Range("H98").Value = 3280.8125
$ws.Range("I98").Value = 4063.182
$ws.Range("K98").Value = 4063.182
$ws.Range("M98").Value = -2565.182
$ws.Range("H112").Value = 1912.1471
$ws.Range("J112").Value = 1828.0344
$ws.Range("L112").Value = 5484.1032
$ws.Range("N112").Value = -7700.1032
$ws.Range("H116").Value = 36545452
$ws.Range("J116").Value = 30325514
$ws.Range("L116").Value = 30325514
$ws.Range("N116").Value = -30332398
$ws.Range("H122").Value = 3280.8125
$ws.Range("I122").Value = 4063.182
$ws.Range("K122").Value = 12189.546
$ws.Range("M122").Value = -9739.545999999998
$ws.Range("H132").Value = 3066.0217
$ws.Range("I132").Value = 2648.5908
$ws.Range("K132").Value = 7945.7724
$ws.Range("M132").Value = -5415.7724
$ws.Range("H137").Value = 3000
$ws.Range("J137").Value = 3000
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -14100

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2521.35
$ws.Range("I74").Value = 2489.5334
$ws.Range("J74").Value = 2616.8
$ws.Range("K74").Value = 2489.5334
$ws.Range("L74").Value = 2616.8
$ws.Range("M74").Value = -1615.5334
$ws.Range("N74").Value = -4364.8
$ws.Range("H77").Value = 2521.35
$ws.Range("I77").Value = 2489.5334
$ws.Range("J77").Value = 2616.8
$ws.Range("K77").Value = 12447.667
$ws.Range("L77").Value = 13084
$ws.Range("M77").Value = -8079.666999999999
$ws.Range("N77").Value = -21820
$ws.Range("H132").Value = 40002336
$ws.Range("I132").Value = 41668976
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 125006928
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -125004398
$ws.Range("N132").Value = -14060
$ws.Range("H139").Value = 89499.164
$ws.Range("J139").Value = 89499.164
$ws.Range("L139").Value = 89499.164
$ws.Range("N139").Value = -99779.164

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3444.5557
$ws.Range("I86").Value = 2714.4285
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 2714.4285
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -1591.4285
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 3444.5557
$ws.Range("I89").Value = 2714.4285
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 13572.1425
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -7956.1425
$ws.Range("N89").Value = -41232
$ws.Range("H94").Value = 2345.818
$ws.Range("I94").Value = 2895.5
$ws.Range("J94").Value = 1500.1538
$ws.Range("K94").Value = 2895.5
$ws.Range("L94").Value = 1500.1538
$ws.Range("M94").Value = -2444.5
$ws.Range("N94").Value = -2402.1538
$ws.Range("H105").Value = 1627.2307
$ws.Range("I105").Value = 1435.6
$ws.Range("K105").Value = 1435.6
$ws.Range("M105").Value = 311.4000000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3068.2778
$ws.Range("I31").Value = 1976.4445
$ws.Range("K31").Value = 1976.4445
$ws.Range("M31").Value = -1681.4445
$ws.Range("H34").Value = 3068.2778
$ws.Range("I34").Value = 1976.4445
$ws.Range("K34").Value = 1976.4445
$ws.Range("M34").Value = -1774.4445
$ws.Range("H102").Value = 45899
$ws.Range("J102").Value = 45899
$ws.Range("L102").Value = 45899
$ws.Range("N102").Value = -50767
$ws.Range("H103").Value = 9016.75
$ws.Range("I103").Value = 9016.75
$ws.Range("K103").Value = 9016.75
$ws.Range("M103").Value = -7844.75
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 999
$ws.Range("I36").Value = 999
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2997
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2828
$ws.Range("H97").Value = 1285.5714
$ws.Range("I97").Value = 1576.25
$ws.Range("J97").Value = 898
$ws.Range("K97").Value = 4728.75
$ws.Range("L97").Value = 2694
$ws.Range("M97").Value = -4232.75
$ws.Range("N97").Value = -3686
$ws.Range("H136").Value = 1948.6
$ws.Range("I136").Value = 1764.5
$ws.Range("J136").Value = 1976.9231
$ws.Range("K136").Value = 5293.5
$ws.Range("L136").Value = 5930.7693
$ws.Range("M136").Value = -193.5
$ws.Range("N136").Value = -16130.7693
$ws.Range("N36").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 28000
$ws.Range("I32").Value = 25000
$ws.Range("J32").Value = 29500
$ws.Range("K32").Value = 25000
$ws.Range("L32").Value = 29500
$ws.Range("M32").Value = -24704
$ws.Range("N32").Value = -30092
$ws.Range("H70").Value = 12581.143
$ws.Range("I70").Value = 12029.333
$ws.Range("J70").Value = 13960.667
$ws.Range("K70").Value = 12029.333
$ws.Range("L70").Value = 13960.667
$ws.Range("M70").Value = -11759.333
$ws.Range("N70").Value = -14500.667
$ws.Range("H73").Value = 12581.143
$ws.Range("I73").Value = 12029.333
$ws.Range("J73").Value = 13960.667
$ws.Range("K73").Value = 12029.333
$ws.Range("L73").Value = 13960.667
$ws.Range("M73").Value = -11093.333
$ws.Range("N73").Value = -15832.667
$ws.Range("H80").Value = 2353.4167
$ws.Range("I80").Value = 1938
$ws.Range("J80").Value = 2491.889
$ws.Range("K80").Value = 1938
$ws.Range("L80").Value = 2491.889
$ws.Range("M80").Value = -940
$ws.Range("N80").Value = -4487.889
$ws.Range("H83").Value = 2353.4167
$ws.Range("I83").Value = 1938
$ws.Range("J83").Value = 2491.889
$ws.Range("K83").Value = 9690
$ws.Range("L83").Value = 12459.445
$ws.Range("M83").Value = -4698
$ws.Range("N83").Value = -22443.445
$ws.Range("H102").Value = 2081.6072
$ws.Range("I102").Value = 1127.5555
$ws.Range("J102").Value = 3798.9
$ws.Range("K102").Value = 1127.5555
$ws.Range("L102").Value = 3798.9
$ws.Range("M102").Value = 494.4445000000001
$ws.Range("N102").Value = -7042.9
$ws.Range("H113").Value = 3617.182
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 3828.9
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 3828.9
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -8168.9
$ws.Range("H132").Value = 4278.231
$ws.Range("I132").Value = 3615.8
$ws.Range("J132").Value = 5181.5454
$ws.Range("K132").Value = 10847.4
$ws.Range("L132").Value = 15544.6362
$ws.Range("M132").Value = -8317.400000000001
$ws.Range("N132").Value = -20604.6362

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1981.3334
$ws.Range("J46").Value = 1981.3334
$ws.Range("L46").Value = 1981.3334
$ws.Range("N46").Value = -2357.3334
$ws.Range("H100").Value = 3232.4285
$ws.Range("I100").Value = 2753.8
$ws.Range("J100").Value = 3498.3333
$ws.Range("K100").Value = 2753.8
$ws.Range("L100").Value = 3498.3333
$ws.Range("M100").Value = -2212.8
$ws.Range("N100").Value = -4580.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 53000
$ws.Range("J75").Value = 53000
$ws.Range("L75").Value = 53000
$ws.Range("N75").Value = -54872
$ws.Range("H78").Value = 53000
$ws.Range("J78").Value = 53000
$ws.Range("L78").Value = 159000
$ws.Range("N78").Value = -168360
$ws.Range("H136").Value = 5568.154
$ws.Range("I136").Value = 2348.25
$ws.Range("K136").Value = 7044.75
$ws.Range("M136").Value = -4494.75
